$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "climate change", 334, 0.06, 1),
    @(2, "vulnerability", 218, 0.04, 1.53),
    @(3, "human", 192, 0.04, 1.74),
    @(4, "sustainability", 168, 0.03, 1.99),
    @(5, "adaptive management", 156, 0.03, 2.14),
    @(6, "disaster management", 152, 0.03, 2.2),
    @(7, "risk assessment", 150, 0.03, 2.23),
    @(8, "natural hazard", 115, 0.02, 2.9),
    @(9, "decision making", 109, 0.02, 3.06),
    @(10, "governance approach", 109, 0.02, 3.06),
    @(11, "urban area", 102, 0.02, 3.27),
    @(12, "sustainable development", 102, 0.02, 3.27),
    @(13, "urban planning", 92, 0.02, 3.63),
    @(14, "livelihood", 85, 0.02, 3.93),
    @(15, "conceptual framework", 81, 0.02, 4.12),
    @(16, "social capital", 78, 0.01, 4.28),
    @(17, "disaster", 77, 0.01, 4.34),
    @(18, "stakeholder", 70, 0.01, 4.77),
    @(19, "perception", 70, 0.01, 4.77),
    @(20, "fishery management", 68, 0.01, 4.91),
    @(21, "local participation", 66, 0.01, 5.06),
    @(22, "covid-19", 65, 0.01, 5.14),
    @(23, "ecosystem resilience", 63, 0.01, 5.3),
    @(24, "flood", 62, 0.01, 5.39),
    @(25, "rural area", 62, 0.01, 5.39),
    @(26, "survey", 61, 0.01, 5.48),
    @(27, "urbanization", 61, 0.01, 5.48),
    @(28, "qualitative analysis", 57, 0.01, 5.86),
    @(29, "flooding", 55, 0.01, 6.07),
    @(30, "strategic approach", 55, 0.01, 6.07),
    @(31, "urban development", 54, 0.01, 6.19),
    @(32, "policy making", 51, 0.01, 6.55),
    @(33, "social network analysis", 50, 0.01, 6.68),
    @(34, "spatiotemporal analysis", 49, 0.01, 6.82),
    @(35, "fishing community", 48, 0.01, 6.96),
    @(36, "adaptation", 44, 0.01, 7.59),
    @(37, "female", 41, 0.01, 8.15),
    @(38, "migration", 41, 0.01, 8.15),
    @(39, "environmental change", 41, 0.01, 8.15),
    @(40, "hazard assessment", 39, 0.01, 8.56),
    @(41, "adult", 38, 0.01, 8.79),
    @(42, "environmental policy", 37, 0.01, 9.03),
    @(43, "participatory approach", 37, 0.01, 9.03),
    @(44, "male", 37, 0.01, 9.03),
    @(45, "poverty", 36, 0.01, 9.28),
    @(46, "hazard management", 36, 0.01, 9.28),
    @(47, "government", 34, 0.01, 9.82),
    @(48, "nature-society relations", 34, 0.01, 9.82),
    @(49, "food security", 33, 0.01, 10.12),
    @(50, "comparative study", 32, 0.01, 10.44),
    @(51, "disaster planning", 32, 0.01, 10.44),
    @(52, "socioeconomic conditions", 32, 0.01, 10.44),
    @(53, "psychology", 32, 0.01, 10.44),
    @(54, "spatial analysis", 31, 0.01, 10.77),
    @(55, "risk perception", 31, 0.01, 10.77),
    @(56, "local government", 31, 0.01, 10.77),
    @(57, "urban population", 30, 0.01, 11.13),
    @(58, "neoliberalism", 30, 0.01, 11.13),
    @(59, "coping strategy", 30, 0.01, 11.13),
    @(60, "innovation", 29, 0.01, 11.52),
    @(61, "rural economy", 29, 0.01, 11.52),
    @(62, "water management", 29, 0.01, 11.52),
    @(63, "rural population", 29, 0.01, 11.52),
    @(64, "developing world", 29, 0.01, 11.52),
    @(65, "natural resource", 28, 0.01, 11.93),
    @(66, "geographic information system", 28, 0.01, 11.93),
    @(67, "household survey", 28, 0.01, 11.93),
    @(68, "rural development", 28, 0.01, 11.93),
    @(69, "empirical analysis", 28, 0.01, 11.93),
    @(70, "policy approach", 27, 0.01, 12.37),
    @(71, "regression analysis", 27, 0.01, 12.37),
    @(72, "flood control", 26, 0, 12.85),
    @(73, "coastal zone", 26, 0, 12.85),
    @(74, "resource management", 25, 0, 13.36),
    @(75, "queensland", 25, 0, 13.36),
    @(76, "fishery", 25, 0, 13.36),
    @(77, "water supply", 25, 0, 13.36),
    @(78, "quantitative analysis", 25, 0, 13.36),
    @(79, "socioeconomic status", 25, 0, 13.36),
    @(80, "informal settlement", 25, 0, 13.36),
    @(81, "drought", 25, 0, 13.36),
    @(82, "community development", 24, 0, 13.92),
    @(83, "socioeconomic impact", 24, 0, 13.92),
    @(84, "institutional framework", 24, 0, 13.92),
    @(85, "economic development", 24, 0, 13.92),
    @(86, "collective action", 23, 0, 14.52),
    @(87, "future prospect", 23, 0, 14.52),
    @(88, "extreme event", 23, 0, 14.52),
    @(89, "learning", 23, 0, 14.52),
    @(90, "social media", 22, 0, 15.18),
    @(91, "numerical model", 22, 0, 15.18),
    @(92, "assessment method", 22, 0, 15.18),
    @(93, "knowledge", 22, 0, 15.18),
    @(94, "tourism development", 21, 0, 15.9),
    @(95, "social impact", 21, 0, 15.9),
    @(96, "investment", 21, 0, 15.9),
    @(97, "education", 21, 0, 15.9),
    @(98, "theoretical study", 21, 0, 15.9),
    @(99, "environmental management", 21, 0, 15.9),
    @(100, "urban growth", 21, 0, 15.9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Host "Updated" $data.Count "rows"
